$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

$beforeLen = $tr.Length
$url = "https://github.com/NadyktoA/CVRP/tree/main"

# Append two blank paragraphs and a third paragraph containing the URL.
$tr.InsertAfter("`r`r`r" + $url)

$afterLen = $tr.Length
$urlLen = $afterLen - $beforeLen - 3

# Grab just the newly inserted URL run and turn it into a hyperlink.
$linkRange = $tr.Characters($beforeLen + 4, $urlLen)
$linkRange.Font.Underline = -1
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $url
